$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.886.42"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.643.69"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5065"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.274"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "1.644.19"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "1.870.31"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5640"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "0.0₅7698"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "25.903.06"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.378"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.943"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.131"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.809"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1239"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.806"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.245"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04957"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.296"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.239"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.576"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.390"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9064"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5556"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "1.133.75"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.552"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01569"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9988"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.520"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8039"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("D45").Value = "1.780.31"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -7.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4284"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.784"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05043"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("E51").Value = "  -0.29%  "
